$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.266.39'
$ws.Range("E2").Value = '  +1.16%  '

$ws.Range("D3").Value = '1.610.15'
$ws.Range("E3").Value = '  +0.60%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = "'213.10"
$ws.Range("E5").Value = '  +0.32%  '

$ws.Range("E6").Value = '  -0.06%  '

$ws.Range("E7").Value = '  +0.67%  '

$ws.Range("D8").Value = "'0.249"
$ws.Range("E8").Value = '  +0.81%  '

$ws.Range("E9").Value = '  +0.31%  '

$ws.Range("E10").Value = '  +1.95%  '

$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("D12").Value = '1.834.98'
$ws.Range("E12").Value = '  +0.66%  '

$ws.Range("D13").Value = '1.622.69'
$ws.Range("E13").Value = '  +1.45%  '

$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("E15").Value = '  +1.12%  '

$ws.Range("D16").Value = '26.273.82'
$ws.Range("E16").Value = '  +1.15%  '

$ws.Range("D17").Value = "'62.12"
$ws.Range("E17").Value = '  +3.10%  '

$ws.Range("E18").Value = '  +0.80%  '

$ws.Range("E19").Value = '  -0.18%  '

$ws.Range("D20").Value = "'201.16"
$ws.Range("E20").Value = '  -0.09%  '

$ws.Range("E21").Value = '  +0.96%  '

$ws.Range("E22").Value = '  +0.78%  '

$ws.Range("E23").Value = '  +0.59%  '

$ws.Range("D24").Value = "'1.90"
$ws.Range("E24").Value = '  +3.96%  '

$ws.Range("D25").Value = "'143.37"
$ws.Range("E25").Value = '  +1.43%  '

$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("E27").Value = '  -1.03%  '

$ws.Range("D28").Value = "'15.23"
$ws.Range("E28").Value = '  +0.66%  '

$ws.Range("E29").Value = '  +2.33%  '

$ws.Range("E30").Value = '  +5.75%  '

$ws.Range("D31").Value = "'1.17"
$ws.Range("E31").Value = '  +0.39%  '

$ws.Range("E32").Value = '  +2.55%  '

$ws.Range("E33").Value = '  -0.41%  '

$ws.Range("E34").Value = '  +1.03%  '

$ws.Range("E35").Value = '  +0.50%  '

$ws.Range("D36").Value = '1.157.43'
$ws.Range("E36").Value = '  +2.74%  '

$ws.Range("D37").Value = "'0.0167"
$ws.Range("E37").Value = '  +0.70%  '

$ws.Range("E38").Value = '  -0.10%  '

$ws.Range("E39").Value = '  +1.20%  '

$ws.Range("D40").Value = "'0.788"
$ws.Range("E40").Value = '  -0.41%  '

$ws.Range("E41").Value = '  +1.13%  '

$ws.Range("D42").Value = "'5.34"
$ws.Range("E42").Value = '  +4.27%  '

$ws.Range("E43").Value = '  +0.24%  '

$ws.Range("D44").Value = '1.747.02'
$ws.Range("E44").Value = '  +0.66%  '

$ws.Range("D45").Value = "'92.78"
$ws.Range("E45").Value = '  -0.24%  '

$ws.Range("E46").Value = '  +13.89%  '

$ws.Range("E47").Value = '  +0.79%  '

$ws.Range("D48").Value = "'53.81"
$ws.Range("E48").Value = '  +0.96%  '

$ws.Range("E49").Value = '  +0.78%  '

$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("E51").Value = '  -0.39%  '
